# Motorcycle-Shop-Dataset.xlsx edit
# - Update the "Sales" (column D) figures for the 2021 rows on the
#   "Latest Data" sheet.
# - Switch the active/selected tab from "Data" to "Copyright".

$wb = $excel.ActiveWorkbook

# --- Update sales figures on the "Latest Data" sheet -----------------
$latest = $wb.Worksheets.Item("Latest Data")

$sales = @{
    2  = 26200
    3  = 21600
    4  = 18700
    5  = 37600
    6  = 31000
    7  = 35800
    8  = 39500
    9  = 9100
    10 = 2200
    11 = 18400
    12 = 11000
    13 = 4100
    14 = 3100
    15 = 3500
    16 = 33900
    17 = 31500
    18 = 6100
    19 = 3800
    20 = 1100
    21 = 3500
    22 = 19000
    23 = 4600
    24 = 1400
    25 = 6700
    26 = 16000
    27 = 3400
    28 = 12800
    29 = 3500
    30 = 31500
}

foreach ($row in $sales.Keys) {
    $latest.Range("D$row").Value = $sales[$row]
}

# --- Switch the active sheet from "Data" to "Copyright" ---------------
$copyright = $wb.Worksheets.Item("Copyright")
$copyright.Activate()
